$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra status/strategy columns (E:H) - leaving only A:D (ID, Email, EmailPassword, MLBPassword)
$ws.Columns("E:H").Delete()

# Add two new accounts rows with updated email addresses
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "faiyam.daft.154@faiyamrahman.com"
$ws.Range("C4").Value = "n/a"
$ws.Range("D4").Value = "beatthestreak1"

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "water.water.354@faiyamrahman.com"
$ws.Range("C5").Value = "n/a"
$ws.Range("D5").Value = "beatthestreak1"
